$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108; this shifts existing rows 108..174 down to 109..175
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new data record
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(108, 3).Value = "Ñuble"
$ws.Cells.Item(108, 4).Value = 44488
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112023
$ws.Cells.Item(108, 7).Value = "Brócoli"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 120
$ws.Cells.Item(108, 11).Value = 750
$ws.Cells.Item(108, 12).Value = 800
$ws.Cells.Item(108, 13).Value = 775
$ws.Cells.Item(108, 14).Value = "$/unidad"
$ws.Cells.Item(108, 15).Value = "Región del Maule"
$ws.Cells.Item(108, 16).Value = 775
$ws.Cells.Item(108, 17).Value = 1
$ws.Cells.Item(108, 18).Value = "Hortaliza"
